# code fixes for scheduling
$wb = $excel.ActiveWorkbook

$wsTestData = $wb.Worksheets.Item("TestData")
$wsWindows  = $wb.Worksheets.Item("Windows")

# --- TestData sheet: update NegotiatedBy / Network values for row 2 ---
# (order matters for shared-string table allocation: E! must be added before Meeri Cunniff)
$wsTestData.Range("E2").Value = "E!"
$wsTestData.Range("D2").Value = "Meeri Cunniff"

# --- Windows sheet: add a new row (TC2 / W1) ---
$wsWindows.Range("A4").Value = "TC2"
$wsWindows.Range("B4").Value = "W1"
$wsWindows.Range("C4").Value = 43137
$wsWindows.Range("D4").Value = 43404
$wsWindows.Range("E4").Value = 2
$wsWindows.Range("F4").Value = 4

# apply the same formatting as row 3 (dates + quote-prefixed numbers), set
# AFTER the values so the paste doesn't get clobbered by value assignment
$wsWindows.Range("A3:F3").Copy()
$wsWindows.Range("A4:F4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Selections / active sheet ---
$wsTestData.Range("D2").Select()
$wsWindows.Activate()
$wsWindows.Range("A1").Select()
